# Adds new column "l1" (W) after existing "c1" (V) column, with 0 values
# for each data row, per commit "git commit calificacion hasta p23".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column W
$ws.Range("W1").Value = "l1"

# Copy the header style from V1 (bold, centered, bordered) to W1
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for the added column (rows 2-8), all zero
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 23).Value = 0
}
